$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.816.55"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "1.894.19"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7976"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3178"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("E9").Value = "  -4.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07043"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08047"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7688"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.69%  "

$ws.Range("D13").Value = "1.894.75"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.312"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "

$ws.Range("D16").Value = "29.824.57"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.957"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.236"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +19.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9988"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "2.144.69"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9986"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1664"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.344"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.060"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.396"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.534"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.00%  "

$ws.Range("E32").Value = "  +3.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05678"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.055"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7403"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.81%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.624"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01911"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.780"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4412"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.809"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8454"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9984"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "1.035.29"
$ws.Range("E46").Value = "  +4.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.874"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.961"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.73%  "

$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").Value = "2.035.61"
$ws.Range("E51").Value = "  -0.71%  "
